$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (exhibitions) ----
$ws1 = $wb.Worksheets.Item("展览")

# Update "want to go" counts (column F) for existing rows
$ws1.Cells.Item(2, 6).Value = 150
$ws1.Cells.Item(4, 6).Value = 438
$ws1.Cells.Item(6, 6).Value = 60
$ws1.Cells.Item(7, 6).Value = 1286
$ws1.Cells.Item(8, 6).Value = 484
$ws1.Cells.Item(10, 6).Value = 253
$ws1.Cells.Item(11, 6).Value = 165
$ws1.Cells.Item(12, 6).Value = 201
$ws1.Cells.Item(14, 6).Value = 4
$ws1.Cells.Item(15, 6).Value = 13
$ws1.Cells.Item(16, 6).Value = 283
$ws1.Cells.Item(17, 6).Value = 48
$ws1.Cells.Item(18, 6).Value = 228
$ws1.Cells.Item(19, 6).Value = 1610
$ws1.Cells.Item(21, 6).Value = 256
$ws1.Cells.Item(22, 6).Value = 106
$ws1.Cells.Item(23, 6).Value = 888
$ws1.Cells.Item(24, 6).Value = 384
$ws1.Cells.Item(26, 6).Value = 903
$ws1.Cells.Item(27, 6).Value = 1184
$ws1.Cells.Item(28, 6).Value = 59
$ws1.Cells.Item(30, 6).Value = 2769
$ws1.Cells.Item(31, 6).Value = 1541
$ws1.Cells.Item(33, 6).Value = 84
$ws1.Cells.Item(34, 6).Value = 561
$ws1.Cells.Item(35, 6).Value = 844
$ws1.Cells.Item(36, 6).Value = 1559
$ws1.Cells.Item(37, 6).Value = 862
$ws1.Cells.Item(38, 6).Value = 1605
$ws1.Cells.Item(39, 6).Value = 184

# Insert the new exhibition row at position 42 (pushes rows 42-46 down to 43-47)
$ws1.Rows.Item(42).Insert()

# Match the index-column formatting used by every other data row
$ws1.Range("A41").Copy()
$ws1.Range("A42").PasteSpecial(-4122)

# Populate the newly inserted row 42
$ws1.Cells.Item(42, 1).Value = 41
$ws1.Cells.Item(42, 2).NumberFormat = "@"
$ws1.Cells.Item(42, 2).Value = "2024-07-20"
$ws1.Cells.Item(42, 3).Value = "杭州·生如夏花国乙only·日夜场"
$ws1.Cells.Item(42, 4).Value = "祥符街道花园岗街181号 格乐利雅婚礼艺术中心(天空之城店)"
$ws1.Cells.Item(42, 5).Value = "2024.07.20 10:00-07.20 22:30"
$ws1.Cells.Item(42, 6).Value = 1
$ws1.Cells.Item(42, 7).Value = 105
$ws1.Cells.Item(42, 8).Value = "https://show.bilibili.com/platform/detail.html?id=85496"
$ws1.Cells.Item(42, 9).Value = "//i1.hdslb.com/bfs/openplatform/202405/Qut2ZdAi1715411977772.jpeg"

# Fix up the sequential index (column A = row number - 1) for rows shifted down by the insert
$ws1.Cells.Item(43, 1).Value = 42
$ws1.Cells.Item(44, 1).Value = 43
$ws1.Cells.Item(45, 1).Value = 44
$ws1.Cells.Item(46, 1).Value = 45
$ws1.Cells.Item(47, 1).Value = 46

# Update "want to go" counts (column F) for rows shifted down by the insert
$ws1.Cells.Item(43, 6).Value = 750
$ws1.Cells.Item(44, 6).Value = 741
$ws1.Cells.Item(45, 6).Value = 951
$ws1.Cells.Item(46, 6).Value = 403
$ws1.Cells.Item(47, 6).Value = 3280

# ---- Sheet "演出" (performances) ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(10, 6).Value = 9
$ws2.Cells.Item(15, 6).Value = 759
$ws2.Cells.Item(21, 6).Value = 8

# ---- Sheet "全部类型" (all types) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(5, 6).Value = 150
$ws4.Cells.Item(6, 6).Value = 438
$ws4.Cells.Item(8, 6).Value = 60
$ws4.Cells.Item(10, 6).Value = 1286
$ws4.Cells.Item(11, 6).Value = 484
$ws4.Cells.Item(13, 6).Value = 253
$ws4.Cells.Item(14, 6).Value = 165
$ws4.Cells.Item(15, 6).Value = 201
$ws4.Cells.Item(17, 6).Value = 13
$ws4.Cells.Item(18, 6).Value = 283
$ws4.Cells.Item(19, 6).Value = 48
$ws4.Cells.Item(20, 6).Value = 228
$ws4.Cells.Item(21, 6).Value = 1610
$ws4.Cells.Item(23, 6).Value = 256
$ws4.Cells.Item(24, 6).Value = 384
$ws4.Cells.Item(25, 6).Value = 9
$ws4.Cells.Item(28, 6).Value = 1184
$ws4.Cells.Item(29, 6).Value = 2769
$ws4.Cells.Item(31, 6).Value = 1541
$ws4.Cells.Item(33, 6).Value = 759
$ws4.Cells.Item(35, 6).Value = 561
$ws4.Cells.Item(36, 6).Value = 844
$ws4.Cells.Item(37, 6).Value = 1559
$ws4.Cells.Item(38, 6).Value = 8
$ws4.Cells.Item(39, 6).Value = 862
$ws4.Cells.Item(40, 6).Value = 1605
$ws4.Cells.Item(42, 6).Value = 750
$ws4.Cells.Item(43, 6).Value = 741
$ws4.Cells.Item(44, 6).Value = 951
$ws4.Cells.Item(45, 6).Value = 403
$ws4.Cells.Item(48, 6).Value = 3280

